# Apply updated values to the "CE" worksheet of the Balance Sheet workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CE")

# Row 4 - Inventory
$ws.Range("B4").Value = 1025000000.0
$ws.Range("C4").Value = 978000000.0
$ws.Range("D4").Value = 1001000000.0
$ws.Range("E4").Value = 1031000000.0
$ws.Range("F4").Value = 1036000000.0

# Row 14 - Accounts Payable
$ws.Range("B14").Value = 906000000.0
$ws.Range("C14").Value = 797000000.0
$ws.Range("D14").Value = 699000000.0
$ws.Range("E14").Value = 599000000.0
$ws.Range("F14").Value = 724000000.0

# Row 22 - Long Term Tax Liability (Deferred)
$ws.Range("B22").Value = 275000000.0
$ws.Range("C22").Value = 250000000.0
$ws.Range("D22").Value = 77000000.0
$ws.Range("E22").Value = 156000000.0
$ws.Range("F22").Value = 167000000.0

# Row 34 - Net Debt
$ws.Range("G34").Value = 3402000000.0

# Row 35 - Total Debt
$ws.Range("G35").Value = 3905000000.0

$wb.Save()
